# Update Name of Algo
# Applies targeted numeric cell updates to Sheet1 matching the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B3").Value = 6.228
$ws.Range("A9").Value = -20.912
$ws.Range("B11").Value = 6.927
$ws.Range("A18").Value = -21.825
$ws.Range("A20").Value = -21.738
$ws.Range("C21").Value = -12.688
